# Updates cryptos price/volume table: refreshed D (Price) and E (Volume(1h))
# values for rows 2-51, plus a Cardano/Dogecoin row swap (rows 8-9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.059.56"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.645.33"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5044"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06456"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").Value = "1.648.86"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "1.867.47"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5469"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "0.0₅7952"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "26.040.37"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.321"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.981"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.012"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.939"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1157"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.755"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05079"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.245"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.275"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.203"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.548"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.350"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8980"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.622"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5669"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "1.154.18"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01573"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.571"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.012"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8198"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "1.779.79"
$ws.Range("E46").Value = "  -1.22%  "
$ws.Range("D47").Value = "0.0₈111"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4543"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.010"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05047"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
